# Tasks.xlsx update: "Minor updates. Continued level building"
#
# - Rows 20-23 (Completed tasks) get a strikethrough font treatment.
# - Rows 28-30 get their borders normalized to match the rest of the
#   table (same style family as rows 16-27), and row 30 (Completed)
#   also gets the strikethrough treatment.
# - A new task row is appended: "Split UI to multiple canvases".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize formatting of the trailing rows (28, 29) to match the
# rest of the task table (copy the cell formatting used by row 24). ---
$ws.Range("A24:E24").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Range("A24:E24").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)

# Row 30 is a "Completed" item - bring it into the same bordered family
# first, then apply strikethrough below with the other completed rows.
$ws.Range("A24:E24").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)

# --- Mark completed tasks with a strikethrough font. ---
$ws.Range("A20:E23").Font.Strikethrough = $true
$ws.Range("A30:E30").Font.Strikethrough = $true

# --- Add the new task row. ---
$ws.Range("A31").Value = "Split UI to multiple canvases"
$ws.Range("B31").Value = "Not Started"
$ws.Range("C31").Value = 2
$ws.Range("D31").Value = "Task"

$ws.Range("A24:E24").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)

# Re-apply the text now that formatting has been copied over from row 24.
$ws.Range("A31").Value = "Split UI to multiple canvases"
$ws.Range("B31").Value = "Not Started"
$ws.Range("C31").Value = 2
$ws.Range("D31").Value = "Task"

# --- Update the active selection / scroll position to match the
# saved view (best effort). ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A25").Select()
